$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-08 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-09 Thursday", 2)

$d.Content.Find.Execute("849÷5=169, 4", $true, $false, $false, $false, $false, $true, 1, $false, "430÷4=107, 2", 2)
$d.Content.Find.Execute("348÷9=38, 6", $true, $false, $false, $false, $false, $true, 1, $false, "389÷7=55, 4", 2)
$d.Content.Find.Execute("945÷5=189, 0", $true, $false, $false, $false, $false, $true, 1, $false, "656÷6=109, 2", 2)
$d.Content.Find.Execute("229÷9=25, 4", $true, $false, $false, $false, $false, $true, 1, $false, "746÷6=124, 2", 2)
$d.Content.Find.Execute("810÷6=135, 0", $true, $false, $false, $false, $false, $true, 1, $false, "951÷2=475, 1", 2)

$d.Content.Find.Execute("566÷6=94, 2", $true, $false, $false, $false, $false, $true, 1, $false, "443÷8=55, 3", 2)
$d.Content.Find.Execute("445÷8=55, 5", $true, $false, $false, $false, $false, $true, 1, $false, "459÷8=57, 3", 2)
$d.Content.Find.Execute("125÷7=17, 6", $true, $false, $false, $false, $false, $true, 1, $false, "479÷3=159, 2", 2)
$d.Content.Find.Execute("563÷4=140, 3", $true, $false, $false, $false, $false, $true, 1, $false, "713÷9=79, 2", 2)
$d.Content.Find.Execute("151÷4=37, 3", $true, $false, $false, $false, $false, $true, 1, $false, "603÷7=86, 1", 2)

$d.Content.Find.Execute("199÷9=22, 1", $true, $false, $false, $false, $false, $true, 1, $false, "822÷3=274, 0", 2)
$d.Content.Find.Execute("226÷9=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "153÷9=17, 0", 2)
$d.Content.Find.Execute("113÷8=14, 1", $true, $false, $false, $false, $false, $true, 1, $false, "272÷7=38, 6", 2)
$d.Content.Find.Execute("977÷2=488, 1", $true, $false, $false, $false, $false, $true, 1, $false, "530÷9=58, 8", 2)
$d.Content.Find.Execute("392÷9=43, 5", $true, $false, $false, $false, $false, $true, 1, $false, "914÷3=304, 2", 2)

$d.Content.Find.Execute("220÷5=44, 0", $true, $false, $false, $false, $false, $true, 1, $false, "877÷7=125, 2", 2)
$d.Content.Find.Execute("963÷8=120, 3", $true, $false, $false, $false, $false, $true, 1, $false, "648÷8=81, 0", 2)
$d.Content.Find.Execute("741÷8=92, 5", $true, $false, $false, $false, $false, $true, 1, $false, "866÷8=108, 2", 2)
$d.Content.Find.Execute("639÷9=71, 0", $true, $false, $false, $false, $false, $true, 1, $false, "219÷2=109, 1", 2)
$d.Content.Find.Execute("494÷4=123, 2", $true, $false, $false, $false, $false, $true, 1, $false, "426÷9=47, 3", 2)

$d.Content.Find.Execute("776÷7=110, 6", $true, $false, $false, $false, $false, $true, 1, $false, "895÷3=298, 1", 2)
$d.Content.Find.Execute("562÷4=140, 2", $true, $false, $false, $false, $false, $true, 1, $false, "738÷5=147, 3", 2)
$d.Content.Find.Execute("119÷9=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "695÷7=99, 2", 2)
$d.Content.Find.Execute("768÷3=256, 0", $true, $false, $false, $false, $false, $true, 1, $false, "252÷9=28, 0", 2)
$d.Content.Find.Execute("105÷5=21, 0", $true, $false, $false, $false, $false, $true, 1, $false, "726÷2=363, 0", 2)

Write-Output "Done"
